# Update SwaadSutra_Daily_2026-01-21.xlsx
# A new order (#28, Vipula Thakkar) comes in and is inserted as the newest
# row at the top of the Daily Orders log (row 2), pushing the existing
# order (#27, Renu) down to row 3. The Summary and Items Breakdown sheets
# are recalculated to include the new order.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet 1: "Daily Orders"
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Daily Orders")

# Insert a new blank row above the current row 2, shifting the existing
# order (Renu, #27) down to row 3.
$ws1.Rows.Item(2).Insert()

# The text columns (everything except the numeric Order ID / Total
# columns) should stay plain text - force text formatting first so
# values like phone numbers / dates aren't auto-converted by Excel's
# type inference when assigned through .Value.
$ws1.Range("B2:F2").NumberFormat = "@"
$ws1.Range("H2:N2").NumberFormat = "@"

$ws1.Range("A2").Value = 28
$ws1.Range("B2").Value = "2026-01-21 10:01"
$ws1.Range("C2").Value = "Vipula Thakkar"
$ws1.Range("D2").Value = "B-903, Kakkad lavida"
$ws1.Range("E2").Value = "8109861246"
$ws1.Range("F2").Value = "Appe Chutney x1"
$ws1.Range("G2").Value = 60
$ws1.Range("H2").Value = "NEW"
$ws1.Range("I2").Value = "PENDING"
$ws1.Range("J2").Value = "2026-01-21"
$ws1.Range("K2").Value = ""
$ws1.Range("L2").Value = ""
$ws1.Range("M2").Value = ""
$ws1.Range("N2").Value = ""

# Note: the sheet's "numberStoredAsText" ignoredError range also grows
# from A1:N2 to A1:N3 whenever Excel re-validates the sheet after this
# edit (the table now spans one more row). That re-validation happens
# automatically as part of normal error checking and isn't something
# that needs to be (or can be) poked directly via the object model.

# ---------------------------------------------------------------
# Sheet 2: "Summary"
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Summary")

$ws2.Range("A2").Value = 2    # Total Orders
$ws2.Range("B2").Value = 2    # New
$ws2.Range("G2").Value = 170  # Total Revenue

# ---------------------------------------------------------------
# Sheet 3: "Items Breakdown"
# ---------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Items Breakdown")

$ws3.Range("B2").Value = 2    # Appe Chutney quantity ordered
$ws3.Range("C2").Value = 120  # Appe Chutney revenue
